$d = $word.ActiveDocument

$replacements = @(
    @{old = "26×45="; new = "50×93="},
    @{old = "57×34="; new = "56×56="},
    @{old = "69×37="; new = "66×69="},
    @{old = "53×88="; new = "27×59="},
    @{old = "36×56="; new = "49×33="},
    @{old = "31×96="; new = "67×17="},
    @{old = "46×85="; new = "41×13="},
    @{old = "92×54="; new = "27×17="},
    @{old = "84×41="; new = "25×83="},
    @{old = "96×95="; new = "31×61="},
    @{old = "46×23="; new = "87×43="},
    @{old = "48×21="; new = "48×74="},
    @{old = "30×49="; new = "20×48="},
    @{old = "62×70="; new = "85×21="},
    @{old = "80×50="; new = "79×29="},
    @{old = "22×86="; new = "11×34="},
    @{old = "37×67="; new = "39×57="},
    @{old = "70×91="; new = "77×62="},
    @{old = "24×93="; new = "37×76="},
    @{old = "96×99="; new = "95×94="},
    @{old = "18×55="; new = "73×66="},
    @{old = "59×89="; new = "15×18="},
    @{old = "42×51="; new = "66×85="},
    @{old = "12×73="; new = "69×47="},
    @{old = "15×41="; new = "23×77="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
